$wb = $excel.ActiveWorkbook

# --- Step 1: rename existing "总计" sheet to "2022-Q1" (keeps its sheetId/rId slot) ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Name = "2022-Q1"

# --- Step 2: add a brand-new sheet right after it and call it "总计" ---
$wsNewTotal = $wb.Worksheets.Add($null, $wsTotal)
$wsNewTotal.Name = "总计"

# Style donors from an existing, untouched quarter sheet (so we reuse the
# workbook's existing header/index cell style instead of inventing a new one).
$headerSample = $wb.Worksheets.Item("2021-Q4").Range("B1")
$indexSample = $wb.Worksheets.Item("2021-Q4").Range("A2")

# ======================================================================
# Step 3: populate "2022-Q1" with the fund-holdings table (replacing the
# old date-summary rows it inherited from "总计").
# ======================================================================
$ws1 = $wb.Worksheets.Item("2022-Q1")
$ws1.Cells.Clear()

$ws1.Range("B1").Value = "基金代码"
$ws1.Range("C1").Value = "基金名称"
$ws1.Range("D1").Value = "基金规模"
$ws1.Range("E1").Value = "股票总仓位"
$ws1.Range("F1").Value = "仓位占比"
$ws1.Range("G1").Value = "持有市值(亿元)"
$ws1.Range("H1").Value = "仓位排名"

$fundRows = @(
    @("001606", "农银汇理工业4.0灵活配置混合", "50.79", "80.16", "4.61", "2.3414", 8),
    @("000336", "农银研究精选混合", "46.79", "77.21", "4.62", "2.1617", 7),
    @("001645", "国泰大健康股票A", "34.81", "90.83", "4.79", "1.6674", 9),
    @("090001", "大成价值增长混合", "18.65", "61.32", "5.82", "1.0854", 2),
    @("020001", "国泰金鹰增长灵活配置混合", "17.38", "91.77", "6.10", "1.0602", 8),
    @("160215", "国泰价值经典灵活配置混合（LOF）", "6.36", "92.86", "5.68", "0.3612", 8),
    @("160919", "大成产业升级股票(LOF)", "3.95", "87.76", "5.09", "0.2011", 4),
    @("011321", "国泰大健康股票C", "3.47", "90.83", "4.79", "0.1662", 9),
    @("006977", "农银汇理海棠三年定期开放混合", "4.64", "63.10", "3.03", "0.1406", 4),
    @("002567", "大成国家安全主题灵活配置混合", "0.34", "52.90", "4.73", "0.0161", 4)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $ws1.Range("A" + $r).Value = $i

    # B:G all look like numbers ("001606", "50.79", "2.3414", ...) but must
    # round-trip as literal text (matches the source sheets for every other
    # quarter). Force text via NumberFormat "@" while writing, then
    # ClearFormats so no stray "s" style attribute lingers on the cell.
    $textCells = $ws1.Range("B" + $r + ":G" + $r)
    $textCells.NumberFormat = "@"
    $ws1.Range("B" + $r).Value = $row[0]
    $ws1.Range("C" + $r).Value = $row[1]
    $ws1.Range("D" + $r).Value = $row[2]
    $ws1.Range("E" + $r).Value = $row[3]
    $ws1.Range("F" + $r).Value = $row[4]
    $ws1.Range("G" + $r).Value = $row[5]
    $textCells.ClearFormats()

    # H (rank) is numeric.
    $ws1.Range("H" + $r).Value = $row[6]
}

# Re-apply the header/index style (bold, centered, thin border) to row 1 and
# column A, matching every sibling quarter sheet - done via copy/paste-format
# so the existing style slot is reused rather than a new one being created.
$headerSample.Copy()
$ws1.Range("B1:H1").PasteSpecial(-4122)
$indexSample.Copy()
$ws1.Range("A2:A11").PasteSpecial(-4122)

# ======================================================================
# Step 4: populate the new "总计" sheet with the date-summary table,
# inserting "2022-Q1" at the top and pushing the older quarters down.
# ======================================================================
$ws2 = $wb.Worksheets.Item("总计")
$ws2.Cells.Clear()

$ws2.Range("B1").Value = "日期"
$ws2.Range("C1").Value = "持有数量(只)"
$ws2.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 10, 9.2),
    @("2021-Q4", 14, 14.71),
    @("2021-Q3", 17, 14.65),
    @("2021-Q2", 26, 26.94),
    @("2021-Q1", 34, 22.02),
    @("2020-Q4", 24, 11.78)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $ws2.Range("A" + $r).Value = $i
    $ws2.Range("B" + $r).Value = $row[0]
    $ws2.Range("C" + $r).Value = $row[1]
    $ws2.Range("D" + $r).Value = $row[2]
}

$headerSample.Copy()
$ws2.Range("B1:D1").PasteSpecial(-4122)
$indexSample.Copy()
$ws2.Range("A2:A7").PasteSpecial(-4122)

$ws1.Range("A1").Select()
$ws2.Range("A1").Select()
